$wb = $excel.ActiveWorkbook

# --- Service Contacts sheet -------------------------------------------------
$wsSvc = $wb.Worksheets.Item("Service Contacts")

# Widen column A and move the selection to D3 (single cell).
$wsSvc.Columns.Item(1).ColumnWidth = 13.666666666666666
$wsSvc.Range("D3").Select()

# --- Practitioners sheet ----------------------------------------------------
$wsPrac = $wb.Worksheets.Item("Practitioners")

# New practitioner row (row 6), matching the existing rows' layout.
$wsPrac.Range("A6").Value = "PHN999:NFP02"
$wsPrac.Range("B6").Value = "P01"
$wsPrac.Range("C6").Value = 8
$wsPrac.Range("D6").Value = 1
$wsPrac.Range("E6").Value = 1973
$wsPrac.Range("F6").Value = 2
$wsPrac.Range("G6").Value = 1
$wsPrac.Range("H6").Value = 1
$wsPrac.Range("I6").Value = "tag1"

# Widen columns A, C and F.
$wsPrac.Columns.Item(1).ColumnWidth = 13.833333333333334
$wsPrac.Columns.Item(3).ColumnWidth = 12.166666666666666
$wsPrac.Columns.Item(6).ColumnWidth = 12

# Select the whole of column G.
$wsPrac.Columns.Item(7).Select()

# Restore the originally active sheet/tab (Metadata) and its selection, so
# switching sheets above doesn't change which tab is active on save.
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("C2").Select()
